$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Давай"
$ws.Range("C2").Value = "Шаги"
$ws.Range("D2").Value = "Ожидаемый результат"

$ws.Range("D3").Select()
